$wb = $excel.ActiveWorkbook

# Overview sheet: E2:F3 status values
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# zh-cn sheet: C2:C3 status values
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C3").Value = "In Translation"
$wsZh.Columns.Item(3).AutoFit() | Out-Null

# de-de sheet: C2:C3 status values
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C3").Value = "In Translation"
$wsDe.Columns.Item(3).AutoFit() | Out-Null
